# "a lot of new" - refresh the article list:
#  - A2/A3 keep their values (266 / 269) but get a style/number-format touch
#  - A4 is corrected from 271 to 270
#  - A5's old value (270, now a duplicate of A4) is cleared out
#  - cursor ends up parked on F8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply "General" number formatting to the article numbers, which is
# what produces the new cellXfs entry (s="2") on A2:A5.
$ws.Range("A2:A5").NumberFormat = "General"

# Fix up the data itself.
$ws.Range("A4").Value = 270
$ws.Range("A5").ClearContents()

# Leave the selection where the editor left it.
$ws.Range("F8").Select() | Out-Null
